# edit.ps1
# Applies the cryptocurrency price/volume update (Coin price + 1h volume%)
# as scraped by the "Updated cryptos list" GitHub Actions workflow.
#
# Columns: D = Price (stored as text), E = Volume(1h)% (stored as text,
# padded with two leading/trailing spaces). Several D values are plain
# decimal numbers (e.g. "317.21", "0.9999") which Excel would otherwise
# auto-convert to numeric cells on assignment; forcing NumberFormat="@"
# before the write (and clearing the format afterwards, to avoid leaving
# a stray style behind) keeps them as literal text, matching the source
# data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.921.12'
$ws.Range("E2").Value = '  -3.40%  '
$ws.Range("D3").Value = '1.857.47'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.21'
$ws.Range("E5").Value = '  -2.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4348'
$ws.Range("E7").Value = '  -5.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3691'
$ws.Range("E8").Value = '  -3.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07481'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.25'
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").Value = '1.851.93'
$ws.Range("E12").Value = '  -3.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.690'
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.430'
$ws.Range("E14").Value = '  -4.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06840'
$ws.Range("E15").Value = '  -2.96%  '
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.45'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009029'
$ws.Range("E18").Value = '  -4.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9994'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.89'
$ws.Range("E20").Value = '  -4.46%  '
$ws.Range("D21").Value = '27.904.95'
$ws.Range("E21").Value = '  -3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.096'
$ws.Range("E22").Value = '  -4.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '2.060.31'
$ws.Range("E24").Value = '  -4.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.008'
$ws.Range("E25").Value = '  -4.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.85'
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.39'
$ws.Range("E27").Value = '  -3.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.365'
$ws.Range("E28").Value = '  -4.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.47'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.725'
$ws.Range("E30").Value = '  -7.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08978'
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8015'
$ws.Range("E32").Value = '  -7.38%  '
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.985'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("E35").Value = '  -6.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9998'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.114'
$ws.Range("E37").Value = '  -4.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05453'
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01975'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.953'
$ws.Range("E40").Value = '  +2.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5240'
$ws.Range("E41").Value = '  -4.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.985'
$ws.Range("E42").Value = '  -5.83%  '
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.771'
$ws.Range("E44").Value = '  -5.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06724'
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4885'
$ws.Range("E46").Value = '  -5.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.62'
$ws.Range("E47").Value = '  -5.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.76'
$ws.Range("E48").Value = '  -4.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.922'
$ws.Range("E49").Value = '  -8.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.677'
$ws.Range("E50").Value = '  -5.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9993'
$ws.Range("E51").Value = '  -0.30%  '

# Strip the temporary text-number-format back off so styling is untouched
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
